# Refresh after 2 March 2025
# - Rename the "Brown Onions" note label to "Brown Onions (Loose)"
# - Append three new weekly price rows (2025-02-03, 2025-02-16, 2025-03-02)
#   to the PRICES sheet, copying the date-column formatting from the row above
# - Update the active selection to M23

$wb = $excel.ActiveWorkbook

$wsPrices = $wb.Worksheets.Item("PRICES")
$wsNotes  = $wb.Worksheets.Item("NOTES")

# --- Update the "Brown Onions" note label ---------------------------------
$wsNotes.Cells.Item(5, 1).Value = "Brown Onions (Loose)"

# --- Append the new price rows ---------------------------------------------
$lastRow = 16
$newRows = @(
    @(45691, 22.09, 2.5, 3.6, 9,  4.2, 5.7, 1.8, 1.659, 75.95, 14.95, 499.95, 49.95, 109.95, 67.5,  46.31),
    @(45704, 23.02, 4.5, 3.9, 9,  4.2, 5.7, 1.8, 1.699, 75.95, 16.95, 469.95, 49.95, 124.95, 71.95, 46.44),
    @(45718, 22.65, 4.5, 2.5, 11, 4.2, 5.7, 1.8, 1.979, 87.95, 18.95, 579.95, 49.95, 134.95, 71.95, 46.43)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $lastRow + 1

    # Copy the date cell's number format (column A) from the row above so the
    # new date keeps the existing date display style.
    $wsPrices.Range("A$lastRow").Copy()
    $wsPrices.Range("A$row").PasteSpecial(-4122)

    $values = $newRows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $wsPrices.Cells.Item($row, $c).Value = $values[$c - 1]
    }

    $lastRow = $row
}

# --- Update the remembered selection on the PRICES sheet -------------------
$wsPrices.Activate() | Out-Null
$wsPrices.Range("M23").Select() | Out-Null
